# Refresh the cryptos price/volume table (GitHub Actions daily update).
# Cells whose new text looks like a plain number are apostrophe-prefixed so
# they are stored as text (matching the original inline-string cells)
# instead of being auto-converted to a numeric value by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.029.00'
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").Value = '1.656.06'
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("E4").Value = '  -0.25%  '
$ws.Range("D5").Value = "'309.92"
$ws.Range("E5").Value = '  +0.72%  '
$ws.Range("E6").Value = '  -0.18%  '
$ws.Range("E7").Value = '  -0.33%  '
$ws.Range("D8").Value = "'0.3851"
$ws.Range("E8").Value = '  +0.54%  '
$ws.Range("D9").Value = "'51.25"
$ws.Range("E9").Value = '  +3.87%  '
$ws.Range("D10").Value = "'1.361"
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("D11").Value = "'0.9991"
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("D12").Value = "'0.08464"
$ws.Range("E12").Value = '  +0.43%  '
$ws.Range("D13").Value = "'24.13"
$ws.Range("E13").Value = '  +2.27%  '
$ws.Range("D14").Value = "'7.137"
$ws.Range("E14").Value = '  +1.48%  '
$ws.Range("D15").Value = "'7.907"
$ws.Range("E15").Value = '  +4.75%  '
$ws.Range("D16").Value = "'0.00001321"
$ws.Range("E16").Value = '  +3.40%  '
$ws.Range("D17").Value = '1.654.31'
$ws.Range("E17").Value = '  +2.19%  '
$ws.Range("D18").Value = "'94.80"
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").Value = "'0.06980"
$ws.Range("E19").Value = '  +0.80%  '
$ws.Range("D20").Value = "'19.85"
$ws.Range("E20").Value = '  -0.65%  '
$ws.Range("D21").Value = "'6.932"
$ws.Range("E21").Value = '  +2.00%  '
$ws.Range("D22").Value = "'0.9993"
$ws.Range("E22").Value = '  -0.17%  '
$ws.Range("D23").Value = "'13.67"
$ws.Range("E23").Value = '  +2.02%  '
$ws.Range("D24").Value = '24.026.50'
$ws.Range("E24").Value = '  +0.82%  '
$ws.Range("D25").Value = "'2.487"
$ws.Range("E25").Value = '  +1.42%  '
$ws.Range("E26").Value = '  +6.66%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").Value = "'151.64"
$ws.Range("E28").Value = '  -3.26%  '
$ws.Range("D29").Value = "'5.453"
$ws.Range("E29").Value = '  +3.33%  '
$ws.Range("D30").Value = "'139.58"
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = "'7.898"
$ws.Range("E31").Value = '  +1.65%  '
$ws.Range("D32").Value = "'2.487"
$ws.Range("E32").Value = '  +0.05%  '
$ws.Range("D33").Value = '1.835.86'
$ws.Range("E33").Value = '  +2.26%  '
$ws.Range("D34").Value = "'1.044"
$ws.Range("E34").Value = '  +7.37%  '
$ws.Range("D35").Value = "'0.08119"
$ws.Range("E35").Value = '  +0.62%  '
$ws.Range("D36").Value = "'0.02985"
$ws.Range("E36").Value = '  +3.86%  '
$ws.Range("D37").Value = "'6.797"
$ws.Range("E37").Value = '  +3.82%  '
$ws.Range("D38").Value = "'10.92"
$ws.Range("E38").Value = '  +6.03%  '
$ws.Range("E39").Value = '  +1.00%  '
$ws.Range("D40").Value = "'0.09168"
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").Value = "'0.7585"
$ws.Range("E41").Value = '  +1.46%  '
$ws.Range("D42").Value = "'13.51"
$ws.Range("E42").Value = '  -0.43%  '
$ws.Range("D43").Value = "'1.427"
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").Value = "'16.43"
$ws.Range("E44").Value = '  +2.39%  '
$ws.Range("D45").Value = "'0.6967"
$ws.Range("E45").Value = '  +1.06%  '
$ws.Range("D46").Value = "'2.463"
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").Value = "'0.9989"
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("D49").Value = "'0.08308"
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").Value = "'134.79"
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  +1.59%  '
